$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: new expense entry for Arduino Nano
$ws.Range("A16").Value = (Get-Date -Year 2018 -Month 5 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D16").Value = "https://icdayroi.com/arduino-nano-v3-0-atmega328p"
$ws.Range("B16").Value = "Arduino nano V3.0 ATmega328P x2"
$ws.Range("C16").Value = 120000

# Row 17: new expense entry for PCB etching materials
$ws.Range("A17").Value = (Get-Date -Year 2018 -Month 5 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B17").Value = "In mạch, board đồng, nước rửa mạch, rào cái vuông x4"
$ws.Range("C17").Value = 50000
$ws.Range("D17").Value = "Bảo tín"

$ws.Range("A16:D17").Font.Name = "Calibri"
$ws.Range("A16:D17").Font.Size = 11

$ws.Range("A16:D17").Borders.LineStyle = 1

$ws.Range("A16").NumberFormat = "m/d/yyyy"
$ws.Range("A17").NumberFormat = "m/d/yyyy"
$ws.Range("A16:A17").HorizontalAlignment = -4108
$ws.Range("A16:A17").VerticalAlignment = -4108

$ws.Range("B16:B17").HorizontalAlignment = -4108
$ws.Range("B16:B17").VerticalAlignment = -4108
$ws.Range("B16").WrapText = $true

$ws.Range("C16:C17").NumberFormat = "#,##0"
$ws.Range("C16:C17").HorizontalAlignment = -4108
$ws.Range("C16:C17").VerticalAlignment = -4108

$ws.Range("D17").HorizontalAlignment = -4108
$ws.Range("D17").VerticalAlignment = -4108

$wb.Save()
